$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Pereyra"
$ws.Range("E8").Value = "CM"
$ws.Range("F8").Value = "Argentina"
$ws.Range("G8").Value = "Serie A TIM"
$ws.Range("H8").Value = "Udinese"
